$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.005.82'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.917.24'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.03'
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4600'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3831'
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07705'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9818'
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.24'
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.938.25'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.695'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.964'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06969'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.13'
$ws.Range("E17").Value = '  -1.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009494'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.68'
$ws.Range("E19").Value = '  -1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.000.42'
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.328'
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.95'
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.084'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.54'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.10'
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.697'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.70'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.872'
$ws.Range("E29").Value = '  +3.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09318'
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8673'
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.111'
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.251'
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.052'
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05712'
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02046'
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.041'
$ws.Range("E39").Value = '  +11.36%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.516'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5519'
$ws.Range("E41").Value = '  -1.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1755'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.397'
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000002834'
$ws.Range("E44").Value = '  -7.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.182'
$ws.Range("E45").Value = '  +5.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5199'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.22'
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06897'
$ws.Range("E48").Value = '  +1.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.784'
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.39'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").Value = '  -0.29%  '
